# Deploying to gh-pages from @ Alvearie/alvearie-fhir-ig -- refresh the
# StructureDefinition "communication-payload-data" export: bump the
# version, refresh the publication date, swap the placeholder
# Contact/ContactDetail rows for real Publisher/Jurisdiction metadata,
# and refresh the root Extension row's Short/Definition text on the
# Elements sheet to match the new Title/Description.

$wb = $excel.ActiveWorkbook

# ---- Metadata sheet -------------------------------------------------
$ws = $wb.Worksheets.Item("Metadata")

# Version: 5.0.0 -> 6.0.0
$ws.Range("B3").Value = "6.0.0"

# Date: refreshed publication timestamp
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was blank; now populated
$ws.Range("B9").Value = "Alvearie Team"

# The two duplicated "Contact" / "No display for ContactDetail" rows are
# replaced: the first becomes "Jurisdiction" / "United States of America"
# and the (now redundant) second row is removed entirely, which shifts
# "Description" and everything below it up by one row.
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"
$ws.Rows.Item(11).Delete()

# ---- Elements sheet ---------------------------------------------------
$ws2 = $wb.Worksheets.Item("Elements")

# Root "Extension" element row: Short / Definition now mirror the
# StructureDefinition's own Title / Description instead of the generic
# "Extension" / "An Extension" placeholders.
$ws2.Range("K2").Value = "Communication Payload Data"
$ws2.Range("L2").Value = "Contains dynamic data elements and vendor / medium specific elements to be included in the message payload"
